$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4674.968
$ws.Range("I15").Value = 4674.968
$ws.Range("K15").Value = 14024.904
$ws.Range("M15").Value = -13855.904

$ws.Range("H86").Value = 2515.3076
$ws.Range("I86").Value = 2714.1428
$ws.Range("J86").Value = 2283.3333
$ws.Range("K86").Value = 2714.1428
$ws.Range("L86").Value = 2283.3333
$ws.Range("M86").Value = -1591.1428
$ws.Range("N86").Value = -4529.3333

$ws.Range("H88").Value = 2767.04
$ws.Range("J88").Value = 3095.15
$ws.Range("L88").Value = 3095.15
$ws.Range("N88").Value = -3907.15

$ws.Range("H89").Value = 2515.3076
$ws.Range("I89").Value = 2714.1428
$ws.Range("J89").Value = 2283.3333
$ws.Range("K89").Value = 13570.714
$ws.Range("L89").Value = 11416.6665
$ws.Range("M89").Value = -7954.714
$ws.Range("N89").Value = -22648.6665

$ws.Range("H91").Value = 2767.04
$ws.Range("J91").Value = 3095.15
$ws.Range("L91").Value = 3095.15
$ws.Range("N91").Value = -5903.15

$ws.Range("H98").Value = 2473.9412
$ws.Range("I98").Value = 2527.3333
$ws.Range("J98").Value = 2268
$ws.Range("K98").Value = 2527.3333
$ws.Range("L98").Value = 2268
$ws.Range("M98").Value = -1029.3333
$ws.Range("N98").Value = -5264

$ws.Range("H122").Value = 2473.9412
$ws.Range("I122").Value = 2527.3333
$ws.Range("J122").Value = 2268
$ws.Range("K122").Value = 7581.999899999999
$ws.Range("L122").Value = 6804
$ws.Range("M122").Value = -5131.999899999999
$ws.Range("N122").Value = -11704

$ws.Range("H138").Value = 5073.8096
$ws.Range("I138").Value = 3088.6365
$ws.Range("J138").Value = 6139.0244
$ws.Range("K138").Value = 9265.9095
$ws.Range("L138").Value = 18417.0732
$ws.Range("M138").Value = -4125.9095
$ws.Range("N138").Value = -28697.0732

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 4094.375
$ws.Range("I88").Value = 2458.8572
$ws.Range("J88").Value = 5366.4443
$ws.Range("K88").Value = 2458.8572
$ws.Range("L88").Value = 5366.4443
$ws.Range("M88").Value = -2052.8572
$ws.Range("N88").Value = -6178.4443

$ws.Range("H91").Value = 4094.375
$ws.Range("I91").Value = 2458.8572
$ws.Range("J91").Value = 5366.4443
$ws.Range("K91").Value = 2458.8572
$ws.Range("L91").Value = 5366.4443
$ws.Range("M91").Value = -1054.8572
$ws.Range("N91").Value = -8174.4443

$ws.Range("H122").Value = 1522.1428
$ws.Range("I122").Value = 1411.5555
$ws.Range("J122").Value = 1721.2
$ws.Range("K122").Value = 4234.666499999999
$ws.Range("L122").Value = 5163.6
$ws.Range("M122").Value = -1784.666499999999
$ws.Range("N122").Value = -10063.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 4501
$ws.Range("I23").Value = 4501
$ws.Range("K23").Value = 4501
$ws.Range("M23").Value = -4218

$ws.Range("H60").Value = 10666.667
$ws.Range("J60").Value = 11000
$ws.Range("L60").Value = 11000
$ws.Range("N60").Value = -12198

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 41668840
$ws.Range("I34").Value = 318
$ws.Range("J34").Value = 52634240
$ws.Range("K34").Value = 954
$ws.Range("L34").Value = 157902720
$ws.Range("M34").Value = -870
$ws.Range("N34").Value = -157902888

$ws.Range("H35").Value = 4000
$ws.Range("J35").Value = 5000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15576

$ws.Range("H49").Value = 2850
$ws.Range("J49").Value = 2850
$ws.Range("L49").Value = 8550
$ws.Range("N49").Value = -8862

$ws.Range("H57").Value = 2799
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 2799
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 8397
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -9515

$ws.Range("H64").Value = 1078400
$ws.Range("J64").Value = 2001885.8
$ws.Range("L64").Value = 6005657.4
$ws.Range("N64").Value = -6006197.4

$ws.Range("H67").Value = 1078400
$ws.Range("J67").Value = 2001885.8
$ws.Range("L67").Value = 6005657.4
$ws.Range("N67").Value = -6007529.4

$ws.Range("H68").Value = 1043.5054
$ws.Range("I68").Value = 715.9423
$ws.Range("K68").Value = 2147.8269
$ws.Range("M68").Value = -1336.8269

$ws.Range("H71").Value = 1043.5054
$ws.Range("I71").Value = 715.9423
$ws.Range("K71").Value = 6443.4807
$ws.Range("M71").Value = -2387.4807

$ws.Range("H108").Value = 3300
$ws.Range("J108").Value = 3300
$ws.Range("L108").Value = 9900
$ws.Range("N108").Value = -15660

$ws.Range("H131").Value = 2086648.4
$ws.Range("I131").Value = 5039.923
$ws.Range("J131").Value = 2859817.2
$ws.Range("K131").Value = 15119.769
$ws.Range("L131").Value = 8579451.600000001
$ws.Range("M131").Value = -10079.769
$ws.Range("N131").Value = -8589531.600000001

$ws.Range("H137").Value = 44412.54
$ws.Range("I137").Value = 3336
$ws.Range("J137").Value = 70085.375
$ws.Range("K137").Value = 10008
$ws.Range("L137").Value = 210256.125
$ws.Range("M137").Value = -4908
$ws.Range("N137").Value = -220456.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 26333.555
$ws.Range("I80").Value = 6502
$ws.Range("J80").Value = 31999.715
$ws.Range("K80").Value = 6502
$ws.Range("L80").Value = 31999.715
$ws.Range("M80").Value = -5504
$ws.Range("N80").Value = -33995.715

$ws.Range("H83").Value = 26333.555
$ws.Range("I83").Value = 6502
$ws.Range("J83").Value = 31999.715
$ws.Range("K83").Value = 32510
$ws.Range("L83").Value = 159998.575
$ws.Range("M83").Value = -27518
$ws.Range("N83").Value = -169982.575

$ws.Range("H102").Value = 1751.0454
$ws.Range("I102").Value = 1725.2858
$ws.Range("J102").Value = 1796.125
$ws.Range("K102").Value = 1725.2858
$ws.Range("L102").Value = 1796.125
$ws.Range("M102").Value = -103.2858000000001
$ws.Range("N102").Value = -5040.125

$ws.Range("H122").Value = 6806.357
$ws.Range("I122").Value = 7480.8184
$ws.Range("K122").Value = 22442.4552
$ws.Range("M122").Value = -19992.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6671.25
$ws.Range("I9").Value = 492.5
$ws.Range("J9").Value = 12850
$ws.Range("K9").Value = 492.5
$ws.Range("L9").Value = 12850
$ws.Range("M9").Value = -268.5
$ws.Range("N9").Value = -13298

$ws.Range("H40").Value = 2066.5
$ws.Range("I40").Value = 2066.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2066.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1930.5
$ws.Range("N40").ClearContents()

$ws.Range("H82").Value = 2023.7
$ws.Range("I82").Value = 1383.3334
$ws.Range("J82").Value = 2984.25
$ws.Range("K82").Value = 1383.3334
$ws.Range("L82").Value = 2984.25
$ws.Range("M82").Value = -1022.3334
$ws.Range("N82").Value = -3706.25

$ws.Range("H85").Value = 2023.7
$ws.Range("I85").Value = 1383.3334
$ws.Range("J85").Value = 2984.25
$ws.Range("K85").Value = 1383.3334
$ws.Range("L85").Value = 2984.25
$ws.Range("M85").Value = -135.3334
$ws.Range("N85").Value = -5480.25

$ws.Range("H122").Value = 26976
$ws.Range("I122").Value = 34968
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 104904
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -102454
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4266.467
$ws.Range("I81").Value = 1166.6666
$ws.Range("J81").Value = 5041.4165
$ws.Range("K81").Value = 2333.3332
$ws.Range("L81").Value = 10082.833
$ws.Range("M81").Value = -1272.3332
$ws.Range("N81").Value = -12204.833

$ws.Range("H84").Value = 4266.467
$ws.Range("I84").Value = 1166.6666
$ws.Range("J84").Value = 5041.4165
$ws.Range("K84").Value = 11666.666
$ws.Range("L84").Value = 50414.165
$ws.Range("M84").Value = -6362.666000000001
$ws.Range("N84").Value = -61022.165

$ws.Range("H96").Value = 20001300
$ws.Range("I96").Value = 33334366
$ws.Range("J96").Value = 1700
$ws.Range("K96").Value = 33334366
$ws.Range("L96").Value = 1700
$ws.Range("M96").Value = -33332993
$ws.Range("N96").Value = -4446

$ws.Range("H136").Value = 904.8889
$ws.Range("I136").Value = 882.8182
$ws.Range("J136").Value = 1002
$ws.Range("K136").Value = 2648.4546
$ws.Range("L136").Value = 3006
$ws.Range("M136").Value = -98.45460000000003
$ws.Range("N136").Value = -8106
